# Update the solver-optimized parameter values in C2:C14 on Sheet1.
# Everything downstream (I:Q, N:Q, S, X, Z columns) is formula-driven and
# will recalculate automatically once these inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -0.35614338000000001
$ws.Range("C3").Value = 0.051760550000000002
$ws.Range("C4").Value = 0.046790730000000003
$ws.Range("C5").Value = 4.26739128
$ws.Range("C6").Value = 1.06854914
$ws.Range("C7").Value = 0.69438610999999995
$ws.Range("C8").Value = 0.33182292000000002
$ws.Range("C9").Value = -21.335530599999998
$ws.Range("C10").Value = -5.6175926499999997
$ws.Range("C11").Value = -20.45828289
$ws.Range("C12").Value = 1.24519253
$ws.Range("C13").Value = 0.46882705000000002
$ws.Range("C14").Value = 19.79543786

# Hide the previously visible helper columns (W:Z) that held the grading
# reference table.
$ws.Range("W1:Z1").EntireColumn.Hidden = $true

# Move the active selection as recorded at save time.
$ws.Range("V17").Select()
